$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 / 7: swap Starting/Ending SoC (%) values ---
$ws.Range("B6").Value = 99
$ws.Range("B7").Value = 16

# --- Row 8: relabel "Total distance covered" ---
$ws.Range("A8").Value = "Total distance covered (km)"

# --- Row 9: relabel WH/KM ---
$ws.Range("A9").Value = "Total energy consumption(WH/KM)"

# --- Row 10: relabel Total SOC consumed ---
$ws.Range("A10").Value = "Total SOC consumed(%)"

# --- Row 12-14: append units to labels ---
$ws.Range("A12").Value = "Peak Power(kW)"
$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"

# --- Row 15: relabel + flip sign of Regenerative Effectiveness ---
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 0.007660860966486528

# --- Row 16 / 17: swap Lowest/Highest Cell Voltage (label + value) ---
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.414
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 3.107

# --- Row 18-20: append units to labels ---
$ws.Range("A18").Value = "Difference in Cell Voltage(V)"
$ws.Range("A19").Value = "Minimum Temperature(C)"
$ws.Range("A20").Value = "Maximum Temperature(C)"

# --- Row 21: relabel + fill in previously-empty value ---
$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 10

# --- Row 22-27: append units to labels ---
$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"

# --- Row 28 / 29: swap lowest/highest cell temp labels (values stay put) ---
$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("A29").Value = "lowest cell temp(C)"

# --- Row 30: append unit ---
$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

# --- Row 31: relabel to Battery Voltage(V), new value ---
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 53

# --- Row 32: relabel to Total energy charged(kWh), new value ---
$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.560140433055555

# --- Row 33: relabel to Electricity consumption units(kW), new value ---
$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 0.000000110836916244356

# --- Row 34: relabel to Idling time percentage, new value ---
$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 8.756917108791146

# --- Row 35: relabel to Time spent in 0-10 km/h, new value ---
$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 9.695361971589938

# --- Row 36: relabel to Time spent in 10-20 km/h, new value ---
$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 3.465685435563923

# --- Row 37: relabel to Time spent in 20-30 km/h, new value ---
$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 5.125791545438987

# --- Row 38: relabel to Time spent in 30-40 km/h, new value ---
$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 21.672656740259

# --- Row 39: relabel to Time spent in 40-50 km/h, new value ---
$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 10.04050430714815

# --- Row 40: relabel to Time spent in 50-60 km/h, new value ---
$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 34.39157966797877

# --- Row 41: relabel to Time spent in 60-70 km/h, new value ---
$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 6.768783159335958

# --- Row 42: relabel to Time spent in 70-80 km/h, new value ---
$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 0

# --- Row 43: brand-new row, Time spent in 80-90 km/h ---
$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
